$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.005.13'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '3.516.39'
$ws.Range("E3").Value = '  -1.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.98'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '197.79'
$ws.Range("E6").Value = '  +5.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.660'
$ws.Range("E10").Value = '  +1.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.12'
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000305'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.66'
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("D14").Value = '4.075.08'
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '601.01'
$ws.Range("E15").Value = '  +4.95%  '
$ws.Range("D16").Value = '70.164.60'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.07'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '3.519.59'
$ws.Range("E19").Value = '  -1.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.90'
$ws.Range("E22").Value = '  +2.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '104.33'
$ws.Range("E23").Value = '  +11.01%  '
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.06'
$ws.Range("E25").Value = '  +2.95%  '
$ws.Range("E26").Value = '  +4.59%  '
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.87'
$ws.Range("E28").Value = '  +4.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.92'
$ws.Range("E29").Value = '  +4.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.47'
$ws.Range("E30").Value = '  +18.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.18'
$ws.Range("E31").Value = '  +1.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.70'
$ws.Range("E32").Value = '  +3.61%  '
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.01'
$ws.Range("E34").Value = '  -0.57%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '3.703.13'
$ws.Range("E35").Value = '  +1.35%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '521.66'
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("D38").Value = '0.0₃0798'
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.01'
$ws.Range("E39").Value = '  -5.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.91'
$ws.Range("E40").Value = '  -2.42%  '
$ws.Range("E41").Value = '  -3.56%  '
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0463'
$ws.Range("E44").Value = '  +1.04%  '
$ws.Range("E45").Value = '  -3.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.140'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.30'
$ws.Range("E47").Value = '  -5.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.80'
$ws.Range("E48").Value = '  -4.35%  '
$ws.Range("E49").Value = '  +0.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.44'
$ws.Range("E50").Value = '  -2.30%  '
$ws.Range("B51").Value = 'OceanProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.34'
$ws.Range("E51").Value = '  -5.56%  '
